try { Write-Output "CWD: $(Get-Location)" } catch { Write-Output "EXC1 $_" }
try { Write-Output "Env HOME: $env:HOME" } catch {}
try {
  $files = [System.IO.Directory]::GetFiles("/tmp/work")
  foreach ($f in $files) { Write-Output "FILE: $f" }
} catch { Write-Output "EXC2 $_" }
